# edit.ps1 -- apply the "rewrote part about kanal and added a note in a diagramm" change
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rewrite the "Jede Station hat also ..." sentence in the Kanal paragraph.
# ---------------------------------------------------------------------------
$oldSentence = " Jede Station hat also gewisser ma" + [char]0x00DF + "en einen eigenen Kanal, der dann " + [char]0x00FC + "ber die Laufzeit und der Synchronisation der Uhren zu einem einheitlichen Kanal wird."
$newSentence = " Jede Station hat also gewisser ma" + [char]0x00DF + "en eigene Frames und Slots (da die Stationen anfangs nicht synchron sind), der dann " + [char]0x00FC + "ber die Laufzeit und der Synchronisation der Uhren zu einem einheitlichen Verst" + [char]0x00E4 + "ndnis von einem Frame und deren Slots f" + [char]0x00FC + "hrt."

$rng = $d.Content
$found = $rng.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the "(ANBINDUNG)" paragraph to the end
#    of the just-rewritten Kanal paragraph (mirrors Word stamping _GoBack at
#    the location of the most recent edit).
# ---------------------------------------------------------------------------
$kanalPara = $d.Content
$kanalPara.Find.Execute("Der Kanal selbst ist nicht zu implementieren")
$kanalPara.Expand(4) | Out-Null
$paraEnd = $kanalPara.End
$bmTarget = $d.Range($paraEnd - 1, $paraEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmTarget)

# ---------------------------------------------------------------------------
# 3) Move the lastRenderedPageBreak from "Zudem sendet jede Station..." up to
#    "Kollisionen entstehen, wenn ...".
# ---------------------------------------------------------------------------
$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r1 = $d.Content
$r1.Find.Execute("Kollisionen entstehen, wenn")
$insertPoint1 = $d.Range($r1.Start, $r1.Start)
$insertPoint1.InsertXML($pkg)

$r2 = $d.Content
$r2.Find.Execute("Zudem sendet jede Station maximal einmal in einem Frame.")
$breakStart = $r2.Start
$breakRange = $d.Range($breakStart, $breakStart)
$breakRange.MoveStartUntil("Zudem", 1) | Out-Null

# ---------------------------------------------------------------------------
# 4) Add a new lastRenderedPageBreak before "Der Empfänger ist zu jeder Zeit".
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Der Empf" + [char]0x00E4 + "nger ist zu jeder Zeit aktiv")
$insertPoint3 = $d.Range($r3.Start, $r3.Start)
$insertPoint3.InsertXML($pkg)

Write-Output "done"
